# Updates the cryptos list values (price & 1h volume change) as produced by
# the scheduled "Updated cryptos list ... with GitHub Actions" run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper that writes a literal text value into a cell without letting Excel's
# automatic type inference turn numeric-looking strings (e.g. "7.41",
# "0.0000252", "1.00") into real numbers - and without leaving a stray
# NumberFormat/style override behind afterwards.
function Set-TextValue {
    param($cell, [string]$val)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# row 2 - Bitcoin
Set-TextValue $ws.Cells.Item(2, 4) "70.647.86"
Set-TextValue $ws.Cells.Item(2, 5) "  -0.67%  "

# row 3 - Ethereum
Set-TextValue $ws.Cells.Item(3, 4) "3.802.77"
Set-TextValue $ws.Cells.Item(3, 5) "  -1.13%  "

# row 4 - TetherUSD
Set-TextValue $ws.Cells.Item(4, 5) "  -0.05%  "

# row 5 - BNB
Set-TextValue $ws.Cells.Item(5, 4) "708.24"
Set-TextValue $ws.Cells.Item(5, 5) "  +1.72%  "

# row 6 - Solana
Set-TextValue $ws.Cells.Item(6, 4) "170.31"
Set-TextValue $ws.Cells.Item(6, 5) "  -1.82%  "

# row 7 - LidoStakedEther
Set-TextValue $ws.Cells.Item(7, 4) "3.802.51"
Set-TextValue $ws.Cells.Item(7, 5) "  -1.10%  "

# row 8 - USDC
Set-TextValue $ws.Cells.Item(8, 5) "  +0.01%  "

# row 9 - XRP
Set-TextValue $ws.Cells.Item(9, 5) "  -0.93%  "

# row 10 - Dogecoin
Set-TextValue $ws.Cells.Item(10, 5) "  -1.55%  "

# row 11 - Toncoin
Set-TextValue $ws.Cells.Item(11, 4) "7.41"
Set-TextValue $ws.Cells.Item(11, 5) "  +1.77%  "

# row 12 - Cardano
Set-TextValue $ws.Cells.Item(12, 5) "  -1.09%  "

# row 13 - ShibaInu
Set-TextValue $ws.Cells.Item(13, 4) "0.0000252"
Set-TextValue $ws.Cells.Item(13, 5) "  -1.89%  "

# row 14 - Avalanche
Set-TextValue $ws.Cells.Item(14, 4) "36.04"
Set-TextValue $ws.Cells.Item(14, 5) "  -1.13%  "

# row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Cells.Item(15, 4) "4.445.19"
Set-TextValue $ws.Cells.Item(15, 5) "  -1.13%  "

# row 16 - WrappedEther
Set-TextValue $ws.Cells.Item(16, 4) "3.778.37"
Set-TextValue $ws.Cells.Item(16, 5) "  -1.83%  "

# row 17 - WrappedBTC
Set-TextValue $ws.Cells.Item(17, 4) "70.751.58"
Set-TextValue $ws.Cells.Item(17, 5) "  -0.64%  "

# row 18 - TRON
Set-TextValue $ws.Cells.Item(18, 5) "  +0.03%  "

# row 19 - was Chainlink, now Polkadot
Set-TextValue $ws.Cells.Item(19, 2) "Polkadot"
Set-TextValue $ws.Cells.Item(19, 3) "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue $ws.Cells.Item(19, 4) "7.12"
Set-TextValue $ws.Cells.Item(19, 5) "  -1.62%  "

# row 20 - was Polkadot, now Chainlink
Set-TextValue $ws.Cells.Item(20, 2) "Chainlink"
Set-TextValue $ws.Cells.Item(20, 3) "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(20, 4) "17.38"
Set-TextValue $ws.Cells.Item(20, 5) "  -1.97%  "

# row 21 - BitcoinCash
Set-TextValue $ws.Cells.Item(21, 4) "494.68"
Set-TextValue $ws.Cells.Item(21, 5) "  +0.33%  "

# row 22 - Uniswap
Set-TextValue $ws.Cells.Item(22, 4) "10.64"
Set-TextValue $ws.Cells.Item(22, 5) "  -4.86%  "

# row 23 - Polygon
Set-TextValue $ws.Cells.Item(23, 4) "0.729"
Set-TextValue $ws.Cells.Item(23, 5) "  +0.96%  "

# row 24 - Litecoin
Set-TextValue $ws.Cells.Item(24, 4) "84.33"
Set-TextValue $ws.Cells.Item(24, 5) "  -0.81%  "

# row 25 - PEPE
Set-TextValue $ws.Cells.Item(25, 5) "  -1.41%  "

# row 26 - InternetComputer(DFINITY)
Set-TextValue $ws.Cells.Item(26, 5) "  -1.81%  "

# row 27 - RenderToken
Set-TextValue $ws.Cells.Item(27, 4) "10.42"
Set-TextValue $ws.Cells.Item(27, 5) "  -1.52%  "

# row 28 - WrappedeETH
Set-TextValue $ws.Cells.Item(28, 4) "3.955.12"
Set-TextValue $ws.Cells.Item(28, 5) "  -1.16%  "

# row 29 - Dai
Set-TextValue $ws.Cells.Item(29, 5) "  +0.04%  "

# row 30 - Fetch.AI
Set-TextValue $ws.Cells.Item(30, 5) "  -4.88%  "

# row 31 - PancakeSwap
Set-TextValue $ws.Cells.Item(31, 4) "3.09"
Set-TextValue $ws.Cells.Item(31, 5) "  -2.68%  "

# row 32 - ImmutableX
Set-TextValue $ws.Cells.Item(32, 5) "  -2.45%  "

# row 33 - NEARProtocol
Set-TextValue $ws.Cells.Item(33, 4) "7.31"
Set-TextValue $ws.Cells.Item(33, 5) "  -4.17%  "

# row 34 - EthereumClassic
Set-TextValue $ws.Cells.Item(34, 4) "29.03"
Set-TextValue $ws.Cells.Item(34, 5) "  -2.17%  "

# row 35 - Kaspa
Set-TextValue $ws.Cells.Item(35, 4) "0.174"
Set-TextValue $ws.Cells.Item(35, 5) "  -3.66%  "

# row 36 - was Aptos, now Binance-PegBSC-USD
Set-TextValue $ws.Cells.Item(36, 2) "Binance-PegBSC-USD"
Set-TextValue $ws.Cells.Item(36, 3) "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Cells.Item(36, 4) "1.00"
Set-TextValue $ws.Cells.Item(36, 5) "  +0.12%  "

# row 37 - was RenzoRestakedETH, now Aptos
Set-TextValue $ws.Cells.Item(37, 2) "Aptos"
Set-TextValue $ws.Cells.Item(37, 3) "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Cells.Item(37, 4) "9.11"
Set-TextValue $ws.Cells.Item(37, 5) "  -1.95%  "

# row 38 - was Binance-PegBSC-USD, now RenzoRestakedETH
Set-TextValue $ws.Cells.Item(38, 2) "RenzoRestakedETH"
Set-TextValue $ws.Cells.Item(38, 3) "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
Set-TextValue $ws.Cells.Item(38, 4) "3.773.96"
Set-TextValue $ws.Cells.Item(38, 5) "  -0.66%  "

# row 39 - Hedera
Set-TextValue $ws.Cells.Item(39, 5) "  -3.34%  "

# row 40 - Mantle
Set-TextValue $ws.Cells.Item(40, 5) "  +1.23%  "

# row 41 - Stacks
Set-TextValue $ws.Cells.Item(41, 4) "2.30"
Set-TextValue $ws.Cells.Item(41, 5) "  -3.21%  "

# row 43 - dogwifhat
Set-TextValue $ws.Cells.Item(43, 5) "  -4.03%  "

# row 45 - FirstDigitalUSD
Set-TextValue $ws.Cells.Item(45, 5) "  +0.12%  "

# row 46 - FLOKI
Set-TextValue $ws.Cells.Item(46, 4) "0.000321"
Set-TextValue $ws.Cells.Item(46, 5) "  +5.12%  "

# row 47 - Monero
Set-TextValue $ws.Cells.Item(47, 4) "164.83"
Set-TextValue $ws.Cells.Item(47, 5) "  +0.95%  "

# row 48 - Bittensor
Set-TextValue $ws.Cells.Item(48, 4) "424.85"
Set-TextValue $ws.Cells.Item(48, 5) "  +1.54%  "

# row 49 - OKB
Set-TextValue $ws.Cells.Item(49, 4) "48.72"
Set-TextValue $ws.Cells.Item(49, 5) "  +0.22%  "

# row 50 - Cosmos
Set-TextValue $ws.Cells.Item(50, 4) "8.61"
Set-TextValue $ws.Cells.Item(50, 5) "  -0.59%  "

# row 51 - ONDO
Set-TextValue $ws.Cells.Item(51, 5) "  -1.94%  "
